$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: F3 changes from a number to text "120"; G3:M3 are new text cells.
# Leading apostrophe forces Excel to store these numeric-looking entries as
# text instead of re-parsing them as numbers (matches real Excel semantics).
$ws.Range("F3").Value = "'120"
$ws.Range("G3").Value = "'120"
$ws.Range("H3").Value = "'13"
$ws.Range("I3").Value = "'0"
$ws.Range("J3").Value = "'0"
$ws.Range("K3").Value = "'0"
$ws.Range("L3").Value = "'0"
$ws.Range("M3").Value = "'0"
# Drop the "quote prefix" cell format the apostrophe entry leaves behind so
# the cells keep their original (default) style.
$ws.Range("F3:M3").Style = "Normal"

# Row 5: F5/G5 keep their numbers but stored as text; H5's value becomes 1;
# I5/K5 become blank text cells; J5/L5/M5 become text "0".
$ws.Range("F5").Value = "'12"
$ws.Range("G5").Value = "'13"
$ws.Range("H5").Value = "'1"
$ws.Range("I5").Value = "'"
$ws.Range("J5").Value = "'0"
$ws.Range("K5").Value = "'"
$ws.Range("L5").Value = "'0"
$ws.Range("M5").Value = "'0"
$ws.Range("F5:M5").Style = "Normal"
